$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.141.17"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.02%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.856.73"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.80%  "

$ws.Range("E4").Value = "  +0.26%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.08"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.28%  "

$ws.Range("E6").Value = "  +0.21%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4698"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.51%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2820"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.87%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06545"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.33%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.17"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.85%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07788"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.32%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "97.00"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.85%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.866.22"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.34%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.089"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.98%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6702"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.01%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "283.25"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.24%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.170.68"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.98%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9999"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.18%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.455"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.38%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.56"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.38%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.106.14"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.40%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.000007228"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.34%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.30%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.142"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.47%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "167.60"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.08%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.299"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.43%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.02"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.88%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.925"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -7.53%  "

$ws.Range("E29").Value = "  -2.79%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09629"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.48%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.402"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.80%  "

$ws.Range("E32").Value = "  -2.77%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.086"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.49%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04676"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.66%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6965"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.05%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.086"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.40%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9993"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.26%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.703"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.43%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01857"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.41%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.281"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.58%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.513"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.85%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "71.95"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.91%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8612"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.23%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.949"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.44%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "104.19"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.65%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9998"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.20%  "

$ws.Range("E47").Value = "  -2.42%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.024.83"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.40%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.207"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.69%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.003"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.65%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "33.78"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.39%  "
